$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Consolidated report: update the "Absent" column (H) for rows where
# attendance was recorded as neither Real/Duplicate/Invalid before.
$ws.Range("H5").Value = 1
$ws.Range("H6").Value = 0
$ws.Range("H15").Value = 1
$ws.Range("H16").Value = 0
